$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 6 (dataset now has 4 data rows instead of 5)
$ws.Rows.Item(6).Delete()

# Update column widths (Excel ColumnWidth = target_stored_width - 1 + 1/12)
$ws.Columns.Item(3).ColumnWidth = 7.083333333333333
$ws.Columns.Item(7).ColumnWidth = 7.083333333333333
$ws.Columns.Item(11).ColumnWidth = 7.083333333333333
$ws.Columns.Item(13).ColumnWidth = 6.083333333333333
$ws.Columns.Item(15).ColumnWidth = 7.083333333333333
$ws.Columns.Item(17).ColumnWidth = 7.083333333333333
$ws.Columns.Item(22).ColumnWidth = 7.083333333333333
$ws.Columns.Item(24).ColumnWidth = 7.083333333333333
$ws.Columns.Item(27).ColumnWidth = 7.083333333333333
$ws.Columns.Item(29).ColumnWidth = 7.083333333333333
$ws.Columns.Item(34).ColumnWidth = 7.083333333333333

# Update data values for rows 2-5
$ws.Range("A2").Value = 45147.50694444445
$ws.Range("B2").Value = 19.697
$ws.Range("C2").Value = 13.276
$ws.Range("D2").Value = 4.055
$ws.Range("E2").Value = 41.787
$ws.Range("F2").Value = 33.752
$ws.Range("G2").Value = 15.501
$ws.Range("H2").Value = 49.487
$ws.Range("I2").Value = 23.851
$ws.Range("J2").Value = 9.981
$ws.Range("K2").Value = 15.136
$ws.Range("L2").Value = 16.473
$ws.Range("M2").Value = 17.17
$ws.Range("N2").Value = 4.948
$ws.Range("O2").Value = 15.414
$ws.Range("P2").Value = 21.544
$ws.Range("Q2").Value = 13.131
$ws.Range("R2").Value = 3.479
$ws.Range("S2").Value = 2.254
$ws.Range("T2").Value = 227.272
$ws.Range("U2").Value = 42.867
$ws.Range("V2").Value = 14.228
$ws.Range("W2").Value = 28.315
$ws.Range("X2").Value = 14.504
$ws.Range("Y2").Value = 2.984
$ws.Range("Z2").Value = 24.958
$ws.Range("AA2").Value = 12.568
$ws.Range("AB2").Value = 11.39
$ws.Range("AC2").Value = 13.349
$ws.Range("AD2").Value = 16.948
$ws.Range("AE2").Value = 3.456
$ws.Range("AF2").Value = 43.884
$ws.Range("AG2").Value = 7.874
$ws.Range("AH2").Value = 17.788
$ws.Range("A3").Value = 45147.51388888889
$ws.Range("B3").Value = 5.765
$ws.Range("C3").Value = 3.621
$ws.Range("D3").Value = 1.453
$ws.Range("E3").Value = 12.255
$ws.Range("F3").Value = 9.608000000000001
$ws.Range("G3").Value = 4.538
$ws.Range("H3").Value = 21.534
$ws.Range("I3").Value = 6.981
$ws.Range("J3").Value = 2.841
$ws.Range("K3").Value = 4.094
$ws.Range("L3").Value = 4.884
$ws.Range("M3").Value = 4.979
$ws.Range("N3").Value = 1.457
$ws.Range("O3").Value = 4.512
$ws.Range("P3").Value = 6.255
$ws.Range("Q3").Value = 4.138
$ws.Range("R3").Value = 1.471
$ws.Range("S3").Value = 0.745
$ws.Range("T3").Value = 61.388
$ws.Range("U3").Value = 12.871
$ws.Range("V3").Value = 4.164
$ws.Range("W3").Value = 8.288
$ws.Range("X3").Value = 4.192
$ws.Range("Y3").Value = 1.139
$ws.Range("Z3").Value = 9.992000000000001
$ws.Range("AA3").Value = 3.678
$ws.Range("AB3").Value = 3.507
$ws.Range("AC3").Value = 4.079
$ws.Range("AD3").Value = 4.916
$ws.Range("AE3").Value = 1.265
$ws.Range("AF3").Value = 19.943
$ws.Range("AG3").Value = 2.198
$ws.Range("AH3").Value = 5.208
$ws.Range("A4").Value = 45147.52083333334
$ws.Range("B4").Value = 7.687
$ws.Range("C4").Value = 5.312
$ws.Range("D4").Value = 1.029
$ws.Range("E4").Value = 16.581
$ws.Range("F4").Value = 13.301
$ws.Range("G4").Value = 6.049
$ws.Range("H4").Value = 23.357
$ws.Range("I4").Value = 9.308
$ws.Range("J4").Value = 3.951
$ws.Range("K4").Value = 5.795
$ws.Range("L4").Value = 6.659
$ws.Range("M4").Value = 6.908
$ws.Range("N4").Value = 1.931
$ws.Range("O4").Value = 6.015
$ws.Range("P4").Value = 8.404
$ws.Range("Q4").Value = 5.303
$ws.Range("R4").Value = 0.981
$ws.Range("S4").Value = 0.5580000000000001
$ws.Range("T4").Value = 84.215
$ws.Range("U4").Value = 16.791
$ws.Range("V4").Value = 5.552
$ws.Range("W4").Value = 11.02
$ws.Range("X4").Value = 5.765
$ws.Range("Y4").Value = 1.174
$ws.Range("Z4").Value = 11.091
$ws.Range("AA4").Value = 4.904
$ws.Range("AB4").Value = 4.501
$ws.Range("AC4").Value = 5.265
$ws.Range("AD4").Value = 6.863
$ws.Range("AE4").Value = 0.773
$ws.Range("AF4").Value = 21.057
$ws.Range("AG4").Value = 3.027
$ws.Range("AH4").Value = 6.942
$ws.Range("A5").Value = 45147.52777777778
$ws.Range("B5").Value = 4.32
$ws.Range("C5").Value = 2.91
$ws.Range("D5").Value = 0.7
$ws.Range("E5").Value = 9.32
$ws.Range("F5").Value = 7.37
$ws.Range("G5").Value = 3.4
$ws.Range("H5").Value = 15.1
$ws.Range("I5").Value = 5.24
$ws.Range("J5").Value = 2.21
$ws.Range("K5").Value = 3.15
$ws.Range("L5").Value = 3.75
$ws.Range("M5").Value = 3.87
$ws.Range("N5").Value = 1.09
$ws.Range("O5").Value = 3.38
$ws.Range("P5").Value = 4.73
$ws.Range("Q5").Value = 3.06
$ws.Range("R5").Value = 0.71
$ws.Range("S5").Value = 0.36
$ws.Range("T5").Value = 44.19
$ws.Range("U5").Value = 9.58
$ws.Range("V5").Value = 3.12
$ws.Range("W5").Value = 6.24
$ws.Range("X5").Value = 3.23
$ws.Range("Y5").Value = 0.74
$ws.Range("Z5").Value = 7.05
$ws.Range("AA5").Value = 2.76
$ws.Range("AB5").Value = 2.57
$ws.Range("AC5").Value = 3
$ws.Range("AD5").Value = 3.85
$ws.Range("AE5").Value = 0.5600000000000001
$ws.Range("AF5").Value = 13.82
$ws.Range("AG5").Value = 1.67
$ws.Range("AH5").Value = 3.91
